$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has columns A (week index), B (rec_yds), C (rec_td),
# D (fumbles), E (fantasy points). The scraper re-run adds two new stat
# columns - height and weight - right after "fumbles" and before
# "fantasy points". Insert two blank columns at E:F, which pushes the
# existing "fantasy points" column (and all of its data) from E to G.
$ws.Columns("E:F").Insert()

# New headers for the inserted columns, plus re-affirm the (now shifted)
# fantasy points header so the row keeps its bold/centered/bordered style.
$ws.Range("E1").Value = "height"
$ws.Range("F1").Value = "weight"
$ws.Range("G1").Value = "fantasy points"

# George Kittle is 6'4" (6.333333333333333 ft) and 250 lbs for every game
# in this 2018 log - fill the new columns down for all 16 game rows.
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 5).Value = 6.333333333333333
    $ws.Cells.Item($r, 6).Value = 250
}
